$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "36.234.16"
$ws.Range("E2").Value = "  -1.43%  "

# Row 3
$ws.Range("D3").Value = "2.039.55"
$ws.Range("E3").Value = "  -2.42%  "

# Row 4
$ws.Range("E4").Value = "  -0.19%  "

# Row 5
$c = $ws.Range("D5")
$c.Value = "'244.79"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.52%  "

# Row 6
$c = $ws.Range("D6")
$c.Value = "'0.662"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +1.83%  "

# Row 7
$ws.Range("B7").Value = "Solana"
$ws.Range("C7").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$c = $ws.Range("D7")
$c.Value = "'56.80"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +2.87%  "

# Row 8
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$c = $ws.Range("D8")
$c.Value = "'1.00"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -0.08%  "

# Row 9
$c = $ws.Range("D9")
$c.Value = "'62.74"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +5.33%  "

# Row 10
$c = $ws.Range("D10")
$c.Value = "'0.365"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -0.62%  "

# Row 11
$c = $ws.Range("D11")
$c.Value = "'0.0744"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -2.80%  "

# Row 12
$ws.Range("E12").Value = "  -3.13%  "

# Row 13
$c = $ws.Range("D13")
$c.Value = "'0.907"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +2.95%  "

# Row 14
$c = $ws.Range("D14")
$c.Value = "'14.16"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -5.38%  "

# Row 15
$ws.Range("D15").Value = "2.329.29"
$ws.Range("E15").Value = "  -2.77%  "

# Row 16
$c = $ws.Range("D16")
$c.Value = "'5.37"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -3.22%  "

# Row 17
$ws.Range("D17").Value = "2.026.52"
$ws.Range("E17").Value = "  -3.09%  "

# Row 18
$ws.Range("B18").Value = "Avalanche"
$ws.Range("C18").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$c = $ws.Range("D18")
$c.Value = "'17.47"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.47%  "

# Row 19
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "36.159.23"
$ws.Range("E19").Value = "  -1.63%  "

# Row 20
$c = $ws.Range("D20")
$c.Value = "'71.40"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -2.38%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0854"
$ws.Range("E21").Value = "  -2.56%  "

# Row 22
$c = $ws.Range("D22")
$c.Value = "'237.06"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.05%  "

# Row 23
$c = $ws.Range("D23")
$c.Value = "'5.17"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -5.43%  "

# Row 24
$c = $ws.Range("D24")
$c.Value = "'1.00"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.44%  "

# Row 25
$c = $ws.Range("D25")
$c.Value = "'2.35"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -2.40%  "

# Row 26
$c = $ws.Range("D26")
$c.Value = "'2.25"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +3.52%  "

# Row 27
$c = $ws.Range("D27")
$c.Value = "'9.26"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -6.34%  "

# Row 28
$c = $ws.Range("D28")
$c.Value = "'164.00"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -2.44%  "

# Row 29
$c = $ws.Range("D29")
$c.Value = "'19.91"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -4.41%  "

# Row 30
$c = $ws.Range("D30")
$c.Value = "'0.121"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -2.19%  "

# Row 31
$c = $ws.Range("D31")
$c.Value = "'1.20"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +0.33%  "

# Row 32
$c = $ws.Range("D32")
$c.Value = "'4.97"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -7.61%  "

# Row 33
$c = $ws.Range("D33")
$c.Value = "'0.0597"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -2.06%  "

# Row 34
$c = $ws.Range("D34")
$c.Value = "'4.40"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -6.88%  "

# Row 35
$c = $ws.Range("D35")
$c.Value = "'0.999"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -0.30%  "

# Row 36
$c = $ws.Range("D36")
$c.Value = "'0.0866"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +3.42%  "

# Row 37
$ws.Range("E37").Value = "  -0.85%  "

# Row 38
$c = $ws.Range("D38")
$c.Value = "'2.20"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -8.31%  "

# Row 39
$c = $ws.Range("D39")
$c.Value = "'5.06"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +3.09%  "

# Row 40
$c = $ws.Range("D40")
$c.Value = "'1.22"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -4.99%  "

# Row 41
$c = $ws.Range("D41")
$c.Value = "'2.87"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -0.25%  "

# Row 42
$c = $ws.Range("D42")
$c.Value = "'0.0214"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -2.89%  "

# Row 43
$c = $ws.Range("D43")
$c.Value = "'1.10"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -6.12%  "

# Row 44
$c = $ws.Range("D44")
$c.Value = "'93.14"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -3.50%  "

# Row 45
$c = $ws.Range("D45")
$c.Value = "'0.0900"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -5.74%  "

# Row 46
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$c = $ws.Range("D46")
$c.Value = "'15.95"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -2.43%  "

# Row 47
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "1.365.39"
$ws.Range("E47").Value = "  +1.74%  "

# Row 48
$c = $ws.Range("D48")
$c.Value = "'7.35"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +3.58%  "

# Row 49
$c = $ws.Range("D49")
$c.Value = "'2.94"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +2.23%  "

# Row 50
$c = $ws.Range("D50")
$c.Value = "'2.25"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -6.94%  "

# Row 51
$c = $ws.Range("D51")
$c.Value = "'45.69"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -0.22%  "
